# Fruta / hortaliza, semanal
#
# The weekly refresh reshuffles the per-record rows (2-25) of the daily
# Cereza sub-dataset: each destination row now carries the Fecha/Variedad/
# Calidad/Volumen/Precios/Unidad/Origen/Precio-Kg/Kg-unidad values that used
# to live on a different row, while rows 8 and 9 stay put. The columns that
# are identical for every record (Mercado ID, Mercado, Region, Codreg, Tipo,
# Producto*, Categoria*) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data which gets reshuffled across rows.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# destRow -> sourceRow : destRow receives the values that used to sit on
# sourceRow (snapshotted from the sheet before any writes happen).
$rowMap = @{
    2  = 14
    3  = 18
    4  = 16
    5  = 17
    6  = 24
    7  = 22
    8  = 8
    9  = 9
    10 = 15
    11 = 10
    12 = 19
    13 = 12
    14 = 13
    15 = 7
    16 = 2
    17 = 5
    18 = 6
    19 = 25
    20 = 4
    21 = 3
    22 = 23
    23 = 11
    24 = 20
    25 = 21
}

# Snapshot every source cell's current value first so writes to one row
# never clobber data another row still needs to read. (Value2 is used
# instead of Value because it round-trips numbers/strings/dates as the
# plain serial/scalar the sheet actually stores.)
$snapshot = @{}
foreach ($row in 2..25) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

foreach ($destRow in 2..25) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcValues[$col]
    }
}
